$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '94.389.63'
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = '3.083.52'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '611.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.11'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.379'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.812'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.43%  '
$ws.Range("D11").Value = '3.079.17'
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").Value = '94.083.38'
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000241'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("D17").Value = '3.652.59'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '3.088.28'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '444.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000190'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.69%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("B28").Value = 'Litecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '84.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = '3.239.41'
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.253'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.180'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.17%  '
$ws.Range("E33").Value = '  -4.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.03'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.10%  '
$ws.Range("E35").Value = '  +32.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.71'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.152'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '486.20'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.440'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("B44").Value = 'MantraDAO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '161.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.674'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("E51").Value = '  +0.12%  '
